$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$arr = New-Object 'object[,]' 13,6
$arr[0,0] = '''2026-01-30'
$arr[0,1] = '12:27:34'
$arr[0,2] = '12:00'
$arr[0,3] = 'Bathroom'
$arr[0,4] = 'No Motion'
$arr[0,5] = 'Inactive'
$arr[1,0] = '''2026-01-30'
$arr[1,1] = '12:27:35'
$arr[1,2] = '12:00'
$arr[1,3] = 'Bathroom'
$arr[1,4] = 'No Motion'
$arr[1,5] = 'Inactive'
$arr[2,0] = '''2026-01-30'
$arr[2,1] = '12:27:40'
$arr[2,2] = '12:00'
$arr[2,3] = 'Bathroom'
$arr[2,4] = 'No Motion'
$arr[2,5] = 'Inactive'
$arr[3,0] = '''2026-01-30'
$arr[3,1] = '12:27:45'
$arr[3,2] = '12:00'
$arr[3,3] = 'Bathroom'
$arr[3,4] = 'No Motion'
$arr[3,5] = 'Inactive'
$arr[4,0] = '''2026-01-30'
$arr[4,1] = '12:27:51'
$arr[4,2] = '12:00'
$arr[4,3] = 'Bathroom'
$arr[4,4] = 'No Motion'
$arr[4,5] = 'Inactive'
$arr[5,0] = '''2026-01-30'
$arr[5,1] = '12:27:55'
$arr[5,2] = '12:00'
$arr[5,3] = 'Bathroom'
$arr[5,4] = 'No Motion'
$arr[5,5] = 'Inactive'
$arr[6,0] = '''2026-01-30'
$arr[6,1] = '12:28:00'
$arr[6,2] = '12:00'
$arr[6,3] = 'Bathroom'
$arr[6,4] = 'No Motion'
$arr[6,5] = 'Inactive'
$arr[7,0] = '''2026-01-30'
$arr[7,1] = '12:28:05'
$arr[7,2] = '12:00'
$arr[7,3] = 'Bathroom'
$arr[7,4] = 'No Motion'
$arr[7,5] = 'Inactive'
$arr[8,0] = '''2026-01-30'
$arr[8,1] = '12:28:10'
$arr[8,2] = '12:00'
$arr[8,3] = 'Bathroom'
$arr[8,4] = 'No Motion'
$arr[8,5] = 'Inactive'
$arr[9,0] = '''2026-01-30'
$arr[9,1] = '12:28:16'
$arr[9,2] = '12:00'
$arr[9,3] = 'Bathroom'
$arr[9,4] = 'No Motion'
$arr[9,5] = 'Inactive'
$arr[10,0] = '''2026-01-30'
$arr[10,1] = '12:28:21'
$arr[10,2] = '12:00'
$arr[10,3] = 'Bathroom'
$arr[10,4] = 'No Motion'
$arr[10,5] = 'Inactive'
$arr[11,0] = '''2026-01-30'
$arr[11,1] = '12:28:26'
$arr[11,2] = '12:00'
$arr[11,3] = 'Bathroom'
$arr[11,4] = 'No Motion'
$arr[11,5] = 'Inactive'
$arr[12,0] = '''2026-01-30'
$arr[12,1] = '12:28:31'
$arr[12,2] = '12:00'
$arr[12,3] = 'Bathroom'
$arr[12,4] = 'No Motion'
$arr[12,5] = 'Inactive'
$ws.Range("A18:F30").Value = $arr

$ws = $wb.Worksheets.Item("Humidity")
$arr = New-Object 'object[,]' 11,6
$arr[0,0] = '''2026-01-30'
$arr[0,1] = '12:27:33'
$arr[0,2] = '12:00'
$arr[0,3] = 'Bathroom'
$arr[0,4] = '''86.8%'
$arr[0,5] = 'Active'
$arr[1,0] = '''2026-01-30'
$arr[1,1] = '12:27:34'
$arr[1,2] = '12:00'
$arr[1,3] = 'Bathroom'
$arr[1,4] = '''87.6%'
$arr[1,5] = 'Active'
$arr[2,0] = '''2026-01-30'
$arr[2,1] = '12:27:41'
$arr[2,2] = '12:00'
$arr[2,3] = 'Bathroom'
$arr[2,4] = '''86.7%'
$arr[2,5] = 'Active'
$arr[3,0] = '''2026-01-30'
$arr[3,1] = '12:27:45'
$arr[3,2] = '12:00'
$arr[3,3] = 'Bathroom'
$arr[3,4] = '''87.5%'
$arr[3,5] = 'Active'
$arr[4,0] = '''2026-01-30'
$arr[4,1] = '12:27:53'
$arr[4,2] = '12:00'
$arr[4,3] = 'Bathroom'
$arr[4,4] = '''87.7%'
$arr[4,5] = 'Active'
$arr[5,0] = '''2026-01-30'
$arr[5,1] = '12:28:09'
$arr[5,2] = '12:00'
$arr[5,3] = 'Bathroom'
$arr[5,4] = '''87.6%'
$arr[5,5] = 'Active'
$arr[6,0] = '''2026-01-30'
$arr[6,1] = '12:28:17'
$arr[6,2] = '12:00'
$arr[6,3] = 'Bathroom'
$arr[6,4] = '''87.6%'
$arr[6,5] = 'Active'
$arr[7,0] = '''2026-01-30'
$arr[7,1] = '12:28:21'
$arr[7,2] = '12:00'
$arr[7,3] = 'Bathroom'
$arr[7,4] = '''86.7%'
$arr[7,5] = 'Active'
$arr[8,0] = '''2026-01-30'
$arr[8,1] = '12:28:25'
$arr[8,2] = '12:00'
$arr[8,3] = 'Bathroom'
$arr[8,4] = '''87.6%'
$arr[8,5] = 'Active'
$arr[9,0] = '''2026-01-30'
$arr[9,1] = '12:28:29'
$arr[9,2] = '12:00'
$arr[9,3] = 'Bathroom'
$arr[9,4] = '''87.6%'
$arr[9,5] = 'Active'
$arr[10,0] = '''2026-01-30'
$arr[10,1] = '12:28:33'
$arr[10,2] = '12:00'
$arr[10,3] = 'Bathroom'
$arr[10,4] = '''87.6%'
$arr[10,5] = 'Active'
$ws.Range("A19:F29").Value = $arr

$ws = $wb.Worksheets.Item("Temperature")
$arr = New-Object 'object[,]' 11,6
$arr[0,0] = '''2026-01-30'
$arr[0,1] = '12:27:34'
$arr[0,2] = '12:00'
$arr[0,3] = 'Bathroom'
$arr[0,4] = '22.8C'
$arr[0,5] = 'Active'
$arr[1,0] = '''2026-01-30'
$arr[1,1] = '12:27:34'
$arr[1,2] = '12:00'
$arr[1,3] = 'Bathroom'
$arr[1,4] = '22.7C'
$arr[1,5] = 'Active'
$arr[2,0] = '''2026-01-30'
$arr[2,1] = '12:27:41'
$arr[2,2] = '12:00'
$arr[2,3] = 'Bathroom'
$arr[2,4] = '22.7C'
$arr[2,5] = 'Active'
$arr[3,0] = '''2026-01-30'
$arr[3,1] = '12:27:45'
$arr[3,2] = '12:00'
$arr[3,3] = 'Bathroom'
$arr[3,4] = '22.6C'
$arr[3,5] = 'Active'
$arr[4,0] = '''2026-01-30'
$arr[4,1] = '12:27:53'
$arr[4,2] = '12:00'
$arr[4,3] = 'Bathroom'
$arr[4,4] = '22.7C'
$arr[4,5] = 'Active'
$arr[5,0] = '''2026-01-30'
$arr[5,1] = '12:28:09'
$arr[5,2] = '12:00'
$arr[5,3] = 'Bathroom'
$arr[5,4] = '22.6C'
$arr[5,5] = 'Active'
$arr[6,0] = '''2026-01-30'
$arr[6,1] = '12:28:17'
$arr[6,2] = '12:00'
$arr[6,3] = 'Bathroom'
$arr[6,4] = '22.7C'
$arr[6,5] = 'Active'
$arr[7,0] = '''2026-01-30'
$arr[7,1] = '12:28:21'
$arr[7,2] = '12:00'
$arr[7,3] = 'Bathroom'
$arr[7,4] = '22.7C'
$arr[7,5] = 'Active'
$arr[8,0] = '''2026-01-30'
$arr[8,1] = '12:28:25'
$arr[8,2] = '12:00'
$arr[8,3] = 'Bathroom'
$arr[8,4] = '22.7C'
$arr[8,5] = 'Active'
$arr[9,0] = '''2026-01-30'
$arr[9,1] = '12:28:29'
$arr[9,2] = '12:00'
$arr[9,3] = 'Bathroom'
$arr[9,4] = '22.7C'
$arr[9,5] = 'Active'
$arr[10,0] = '''2026-01-30'
$arr[10,1] = '12:28:33'
$arr[10,2] = '12:00'
$arr[10,3] = 'Bathroom'
$arr[10,4] = '22.8C'
$arr[10,5] = 'Active'
$ws.Range("A19:F29").Value = $arr

$ws = $wb.Worksheets.Item("mmWave")
$arr = New-Object 'object[,]' 7,6
$arr[0,0] = '''2026-01-30'
$arr[0,1] = '12:27:33'
$arr[0,2] = '12:00'
$arr[0,3] = 'Living Room'
$arr[0,4] = 'PRESENCE_DETECTED'
$arr[0,5] = 'Active'
$arr[1,0] = '''2026-01-30'
$arr[1,1] = '12:27:37'
$arr[1,2] = '12:00'
$arr[1,3] = 'Living Room'
$arr[1,4] = 'PRESENCE_DETECTED'
$arr[1,5] = 'Active'
$arr[2,0] = '''2026-01-30'
$arr[2,1] = '12:27:47'
$arr[2,2] = '12:00'
$arr[2,3] = 'Living Room'
$arr[2,4] = 'PRESENCE_DETECTED'
$arr[2,5] = 'Active'
$arr[3,0] = '''2026-01-30'
$arr[3,1] = '12:27:58'
$arr[3,2] = '12:00'
$arr[3,3] = 'Living Room'
$arr[3,4] = 'PRESENCE_DETECTED'
$arr[3,5] = 'Active'
$arr[4,0] = '''2026-01-30'
$arr[4,1] = '12:28:08'
$arr[4,2] = '12:00'
$arr[4,3] = 'Living Room'
$arr[4,4] = 'PRESENCE_DETECTED'
$arr[4,5] = 'Active'
$arr[5,0] = '''2026-01-30'
$arr[5,1] = '12:28:19'
$arr[5,2] = '12:00'
$arr[5,3] = 'Living Room'
$arr[5,4] = 'PRESENCE_DETECTED'
$arr[5,5] = 'Active'
$arr[6,0] = '''2026-01-30'
$arr[6,1] = '12:28:30'
$arr[6,2] = '12:00'
$arr[6,3] = 'Living Room'
$arr[6,4] = 'PRESENCE_DETECTED'
$arr[6,5] = 'Active'
$ws.Range("A13:F19").Value = $arr
